$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old email address used for the "Invalid Login" row
$ws.Range("A3").Value = "bakomam596@skrak.com"

# New "userName" column (C) with header + two data rows
$ws.Range("C1").Value = "userName"
$ws.Range("C2").Value = "Some Name"
$ws.Range("C3").Value = "Jatin Sharma"

# Turn that new email address into a mailto hyperlink
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:bakomam596@skrak.com")

# Resize columns to fit the new layout
$ws.Columns.Item(1).ColumnWidth = 27.21875
$ws.Columns.Item(2).ColumnWidth = 12.109375
$ws.Columns.Item(3).ColumnWidth = 16.88671875

# Print/page layout tweak
$ws.PageSetup.Orientation = 1

# Leave selection where the author left it
$ws.Range("D6").Select()
